# Revert "My 25 added"
# Removes the First name / Surname values that had been added into
# columns C and D for rows 21-45 of the "Patient details" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Patient details")

$ws.Range("C21:D45").ClearContents()

# Restore the previously-selected cell / view state.
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 100
$ws.Range("I10").Select() | Out-Null
